$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF"), matching the existing header style (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data for I2:I59 and J2:J59 (row 2 maps to index 0)
$iVals = @(9,9,9,9,9,9,8,9,9,9,9,9,9,9,8,9,8,9,9,9,8,9,8,9,10,9,9,9,9,9,9,9,8,9,8,9,9,9,9,9,9,9,9,9,9,8,9,9,9,9,9,9,9,9,7,7,5,6)
$jVals = @(9,9,9,9,9,9,8,9,9,9,9,9,9,9,9,9,9,9,9,9,8,9,9,9,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,9,9,9,9,9,9,9,9,7,7,5,6)

for ($k = 0; $k -lt $iVals.Count; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$k]
    $ws.Cells.Item($row, 10).Value = $jVals[$k]
}
